$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShape($range, $newName) {
    $shp = $range.InlineShapes.Item(1)
    $tmp = $shp.ConvertToShape()
    $tmp.Name = $newName
    $tmp.ConvertToInlineShape() | Out-Null
}

# Footer (default / footer2.xml) - Pearson logo: image1.png -> image2.png
Rename-InlineShape $sec.Footers.Item(1).Range "image2.png"

# Footer (first page / footer1.xml) - Pearson logo: image1.png -> image2.png
Rename-InlineShape $sec.Footers.Item(2).Range "image2.png"

# Header (first page / header1.xml) - BTec logo: image2.jpg -> image1.jpg
Rename-InlineShape $sec.Headers.Item(2).Range "image1.jpg"
